# Apply green highlighting to four specific "selected tool" related
# paragraphs in the document, per the commit message
# "green effect added on selected tools".
#
# wdBrightGreen (WdColorIndex) = 4  ->  OOXML <w:highlight w:val="green"/>

$d = $word.ActiveDocument
$wdBrightGreen = 4

# ---------------------------------------------------------------------
# 1) "Even when re-editing and previously edited image, ..."
# 2) "Green outline for the selected button and the selected image (Nice to have)"
# 3) "with a prompt to go back to the screen"
# These three paragraphs only need highlight applied to the whole paragraph
# (including the paragraph mark) - no other structural changes are needed.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if (($text -like "*re-editing and previously edited image*") -or `
        ($text -like "*Green outline for the selected button*") -or `
        ($text -like "*with a prompt to go back to the screen*")) {
        $p.Range.Font.HighlightColorIndex = $wdBrightGreen
    }
}

# ---------------------------------------------------------------------
# 4) "Add thumbs folder (previous changes doc I sent) – (Very nice to have)"
# This paragraph also needs its trailing four runs
# ( " ", "-", " ", "(Very nice to have)" ) merged into a single run
# reading " - (Very nice to have)" before the highlight is applied.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "*Add thumbs folder*") {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $firstRunText = "Add thumbs folder (previous changes doc I sent)"
        $splitAt = $pStart + $firstRunText.Length

        # Temporarily toggle Bold on the first run so that the upcoming
        # replace of the trailing text does not get merged back into it.
        $firstRunRange = $d.Range($pStart, $splitAt)
        $firstRunRange.Font.Bold = 1

        $tailRange = $d.Range($splitAt, $pEnd - 1)
        $tailRange.Find.Execute(" – (Very nice to have)", $false, $false, $false, $false, $false, $true, 1, $false, " – (Very nice to have)", 2)

        # Restore the first run's original (non-bold) formatting.
        $firstRunRange2 = $d.Range($pStart, $splitAt)
        $firstRunRange2.Font.Bold = 0

        # Now apply the green highlight across the whole paragraph.
        $p.Range.Font.HighlightColorIndex = $wdBrightGreen
    }
}
